$d = $word.ActiveDocument

$pairs = @(
    @("402×6=2412", "858×3=2574"),
    @("884×5=4420", "374×2=748"),
    @("510×2=1020", "956×3=2868"),
    @("177×6=1062", "418×9=3762"),
    @("849×2=1698", "214×9=1926"),
    @("765×9=6885", "806×8=6448"),
    @("901×6=5406", "847×7=5929"),
    @("850×7=5950", "953×3=2859"),
    @("675×9=6075", "922×4=3688"),
    @("291×5=1455", "502×8=4016"),
    @("269×4=1076", "573×7=4011"),
    @("414×8=3312", "758×6=4548"),
    @("838×5=4190", "414×5=2070"),
    @("341×6=2046", "616×4=2464"),
    @("696×8=5568", "693×5=3465"),
    @("680×4=2720", "284×4=1136"),
    @("837×2=1674", "853×2=1706"),
    @("431×4=1724", "233×3=699"),
    @("180×3=540", "468×7=3276"),
    @("548×3=1644", "847×7=5929"),
    @("667×6=4002", "426×2=852"),
    @("142×3=426", "208×4=832"),
    @("192×3=576", "514×2=1028"),
    @("142×4=568", "838×3=2514"),
    @("111×2=222", "270×2=540")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
